# Add a "PRESUPUESTO" column (G) to the "VENTA MENSUAL" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# Copy formatting from the neighboring column (F) into the new column (G)
# for each existing row, then set the new values/text.
$ws.Range("F1").Copy($ws.Range("G1"))
$ws.Range("G1").Value = "PRESUPUESTO"

$ws.Range("F2").Copy($ws.Range("G2"))
$ws.Range("G2").Value = 0

$ws.Range("F3").Copy($ws.Range("G3"))
$ws.Range("G3").Value = 0

$ws.Range("F4").Copy($ws.Range("G4"))
$ws.Range("G4").Value = 0

# Column width for the new column (OOXML stores width ~0.83 wider than the
# COM ColumnWidth value for this font/runtime, so back that offset out to
# land on the target stored width of 17).
$ws.Columns.Item(7).ColumnWidth = 16.17
